# Update zhongshu_wangge.xlsx per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Block 1: 创业板50 (rows 7-9) ---
$ws.Range("D8").Value = "1.290/1.305"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = $null

# --- Block 2: 300ETF (rows 16-18) ---
$ws.Range("D17").Value = "4.918/5.074"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = $null

# --- Block 3: 科创50 (rows 25-27) ---
$ws.Range("D26").Value = "1.486/1.524"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = $null

# --- Block 4: 证券ETF (rows 34-36) ---
$ws.Range("B35").Value = "1.211/1.250"
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = $null

# --- Block 5: 生物医药 (rows 43-45) ---
$ws.Range("C45").Value = 1

# --- Block 6: 银行ETF (rows 52-54) ---
$ws.Range("C54").Value = 1

# --- Update the saved view (scroll position / active selection) ---
$ws.Range("D38").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
